$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Save" in H1, matching the bold/bordered header style used by
# the other header cells (copy format from G1, same as B1:G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "Save"

# Fill H2:H38 with a flag derived from the "sum" column G: 1 when the sum
# exceeds 10, otherwise 0.
for ($r = 2; $r -le 38; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -gt 10) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
